# Auto-generated edit script applying scheduled market-price refresh
# values to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Source data has no formulas (static <v> numbers refreshed by an external
# market-data runner), so this just rewrites the affected cells directly.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: ALC ---
$ws = $wb.Worksheets.Item(1)
# Row 107
$ws.Cells.Item(107, 8).Value = 884.9666999999999  # H107: 876.0323 -> 884.9666999999999
$ws.Cells.Item(107, 9).Value = 861.3077  # I107: 851.9259 -> 861.3077
$ws.Cells.Item(107, 11).Value = 861.3077  # K107: 851.9259 -> 861.3077
$ws.Cells.Item(107, 13).Value = 1058.6923  # M107: 1068.0741 -> 1058.6923
# Row 116
$ws.Cells.Item(116, 8).Value = 10751.5  # H116: 7985.1665 -> 10751.5
$ws.Cells.Item(116, 9).Value = 9500  # I116: 6999.6665 -> 9500
$ws.Cells.Item(116, 10).Value = 12003  # J116: 8970.666999999999 -> 12003
$ws.Cells.Item(116, 11).Value = 9500  # K116: 6999.6665 -> 9500
$ws.Cells.Item(116, 12).Value = 12003  # L116: 8970.666999999999 -> 12003
$ws.Cells.Item(116, 13).Value = -6058  # M116: -3557.6665 -> -6058
$ws.Cells.Item(116, 14).Value = -18887  # N116: -15854.667 -> -18887
# Row 137
$ws.Cells.Item(137, 8).Value = 4263.7856  # H137: 4263.857 -> 4263.7856
$ws.Cells.Item(137, 9).Value = 5528  # I137: 5528.143 -> 5528
$ws.Cells.Item(137, 11).Value = 16584  # K137: 16584.429 -> 16584
$ws.Cells.Item(137, 13).Value = -14034  # M137: -14034.429 -> -14034
# Row 138
$ws.Cells.Item(138, 8).Value = 5418.323  # H138: 5454.54 -> 5418.323
$ws.Cells.Item(138, 10).Value = 5562.7046  # J138: 5601.7754 -> 5562.7046
$ws.Cells.Item(138, 12).Value = 16688.1138  # L138: 16805.3262 -> 16688.1138
$ws.Cells.Item(138, 14).Value = -26968.1138  # N138: -27085.3262 -> -26968.1138
# Row 141
$ws.Cells.Item(141, 8).Value = 779.1111  # H141: 748.6316 -> 779.1111
$ws.Cells.Item(141, 9).Value = 779.1111  # I141: 748.6316 -> 779.1111
$ws.Cells.Item(141, 11).Value = 2337.3333  # K141: 2245.8948 -> 2337.3333
$ws.Cells.Item(141, 13).Value = 2842.6667  # M141: 2934.1052 -> 2842.6667

# --- Sheet 2: ARM ---
$ws = $wb.Worksheets.Item(2)
# Row 45
$ws.Cells.Item(45, 8).Value = 11661.5  # H45: 14283.125 -> 11661.5
$ws.Cells.Item(45, 9).Value = 18186  # I45: 21553.2 -> 18186
$ws.Cells.Item(45, 10).Value = 1874.75  # J45: 2166.3333 -> 1874.75
$ws.Cells.Item(45, 11).Value = 18186  # K45: 21553.2 -> 18186
$ws.Cells.Item(45, 12).Value = 1874.75  # L45: 2166.3333 -> 1874.75
$ws.Cells.Item(45, 13).Value = -17809  # M45: -21176.2 -> -17809
$ws.Cells.Item(45, 14).Value = -2628.75  # N45: -2920.3333 -> -2628.75
# Row 61
$ws.Cells.Item(61, 8).Value = 300854.53  # H61: 324813.9 -> 300854.53
$ws.Cells.Item(61, 9).Value = 3419.2104  # I61: 3661.2354 -> 3419.2104
$ws.Cells.Item(61, 11).Value = 3419.2104  # K61: 3661.2354 -> 3419.2104
$ws.Cells.Item(61, 13).Value = -3207.2104  # M61: -3449.2354 -> -3207.2104
# Row 74
$ws.Cells.Item(74, 8).Value = 67799.766  # H74: 54989.117 -> 67799.766
$ws.Cells.Item(74, 9).Value = 86830.92  # I74: 75401.60000000001 -> 86830.92
$ws.Cells.Item(74, 10).Value = 36874.125  # J74: 27153.908 -> 36874.125
$ws.Cells.Item(74, 11).Value = 86830.92  # K74: 75401.60000000001 -> 86830.92
$ws.Cells.Item(74, 12).Value = 36874.125  # L74: 27153.908 -> 36874.125
$ws.Cells.Item(74, 13).Value = -85956.92  # M74: -74527.60000000001 -> -85956.92
$ws.Cells.Item(74, 14).Value = -38622.125  # N74: -28901.908 -> -38622.125
# Row 77
$ws.Cells.Item(77, 8).Value = 67799.766  # H77: 54989.117 -> 67799.766
$ws.Cells.Item(77, 9).Value = 86830.92  # I77: 75401.60000000001 -> 86830.92
$ws.Cells.Item(77, 10).Value = 36874.125  # J77: 27153.908 -> 36874.125
$ws.Cells.Item(77, 11).Value = 434154.6  # K77: 377008 -> 434154.6
$ws.Cells.Item(77, 12).Value = 184370.625  # L77: 135769.54 -> 184370.625
$ws.Cells.Item(77, 13).Value = -429786.6  # M77: -372640 -> -429786.6
$ws.Cells.Item(77, 14).Value = -193106.625  # N77: -144505.54 -> -193106.625
# Row 136
$ws.Cells.Item(136, 8).Value = 300854.53  # H136: 324813.9 -> 300854.53
$ws.Cells.Item(136, 9).Value = 3419.2104  # I136: 3661.2354 -> 3419.2104
$ws.Cells.Item(136, 11).Value = 10257.6312  # K136: 10983.7062 -> 10257.6312
$ws.Cells.Item(136, 13).Value = -7707.6312  # M136: -8433.706200000001 -> -7707.6312

# --- Sheet 3: BSM ---
$ws = $wb.Worksheets.Item(3)
# Row 100
$ws.Cells.Item(100, 8).Value = 25425.334  # H100: 24435.857 -> 25425.334
$ws.Cells.Item(100, 10).Value = 25425.334  # J100: 24435.857 -> 25425.334
$ws.Cells.Item(100, 12).Value = 25425.334  # L100: 24435.857 -> 25425.334
$ws.Cells.Item(100, 14).Value = -27589.334  # N100: -26599.857 -> -27589.334
# Row 103
$ws.Cells.Item(103, 8).Value = 25000  # H103: 24999.5 -> 25000
$ws.Cells.Item(103, 10).Value = 25000  # J103: 24999.5 -> 25000
$ws.Cells.Item(103, 12).Value = 25000  # L103: 24999.5 -> 25000
$ws.Cells.Item(103, 14).Value = -27344  # N103: -27343.5 -> -27344
# Row 134
$ws.Cells.Item(134, 8).Value = 2926.3333  # H134: 2560.7585 -> 2926.3333
$ws.Cells.Item(134, 9).Value = 1840.1875  # I134: 1593.9524 -> 1840.1875
$ws.Cells.Item(134, 11).Value = 5520.5625  # K134: 4781.857199999999 -> 5520.5625
$ws.Cells.Item(134, 13).Value = -2985.5625  # M134: -2246.857199999999 -> -2985.5625

# --- Sheet 4: CRP ---
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Cells.Item(31, 8).Value = 3717.5789  # H31: 2778.5518 -> 3717.5789
$ws.Cells.Item(31, 9).Value = 3241.7693  # I31: 2264.652 -> 3241.7693
$ws.Cells.Item(31, 11).Value = 3241.7693  # K31: 2264.652 -> 3241.7693
$ws.Cells.Item(31, 13).Value = -2946.7693  # M31: -1969.652 -> -2946.7693
# Row 34
$ws.Cells.Item(34, 8).Value = 3717.5789  # H34: 2778.5518 -> 3717.5789
$ws.Cells.Item(34, 9).Value = 3241.7693  # I34: 2264.652 -> 3241.7693
$ws.Cells.Item(34, 11).Value = 3241.7693  # K34: 2264.652 -> 3241.7693
$ws.Cells.Item(34, 13).Value = -3039.7693  # M34: -2062.652 -> -3039.7693
# Row 42
$ws.Cells.Item(42, 8).Value = 6950  # H42: 7450 -> 6950
$ws.Cells.Item(42, 9).Value = 6950  # I42: 7900 -> 6950
$ws.Cells.Item(42, 10).Value = 0  # J42: 7000 -> 0
$ws.Cells.Item(42, 11).Value = 6950  # K42: 7900 -> 6950
$ws.Cells.Item(42, 12).Value = 0  # L42: 7000 -> 0
$ws.Cells.Item(42, 13).Value = -6357  # M42: -7307 -> -6357
$ws.Cells.Item(42, 14).ClearContents()  # N42: remove (was -8186)

# --- Sheet 5: CUL ---
$ws = $wb.Worksheets.Item(5)
# Row 110
$ws.Cells.Item(110, 8).Value = 29666.334  # H110: 31000 -> 29666.334
$ws.Cells.Item(110, 9).Value = 29666.334  # I110: 31000 -> 29666.334
$ws.Cells.Item(110, 11).Value = 88999.00199999999  # K110: 93000 -> 88999.00199999999
$ws.Cells.Item(110, 13).Value = -84909.00199999999  # M110: -88910 -> -84909.00199999999

# --- Sheet 6: GSM ---
$ws = $wb.Worksheets.Item(6)
# Row 43
$ws.Cells.Item(43, 8).Value = 466605.7  # H43: 555060.25 -> 466605.7
$ws.Cells.Item(43, 9).Value = 502156.25  # I43: 669518.3 -> 502156.25
$ws.Cells.Item(43, 11).Value = 502156.25  # K43: 669518.3 -> 502156.25
$ws.Cells.Item(43, 13).Value = -502005.25  # M43: -669367.3 -> -502005.25
# Row 46
$ws.Cells.Item(46, 8).Value = 0  # H46: 50000 -> 0
$ws.Cells.Item(46, 10).Value = 0  # J46: 50000 -> 0
$ws.Cells.Item(46, 12).Value = 0  # L46: 50000 -> 0
$ws.Cells.Item(46, 14).ClearContents()  # N46: remove (was -50312)
# Row 52
$ws.Cells.Item(52, 8).Value = 23998.5  # H52: 25999 -> 23998.5
$ws.Cells.Item(52, 10).Value = 0  # J52: 30000 -> 0
$ws.Cells.Item(52, 12).Value = 0  # L52: 30000 -> 0
$ws.Cells.Item(52, 14).ClearContents()  # N52: remove (was -30518)
# Row 57
$ws.Cells.Item(57, 8).Value = 9999.6  # H57: 6666.6665 -> 9999.6
$ws.Cells.Item(57, 10).Value = 14999  # J57: 0 -> 14999
$ws.Cells.Item(57, 12).Value = 14999  # L57: 0 -> 14999
$ws.Cells.Item(57, 14).Value = -16639  # N57: add (new -16639)
# Row 113
$ws.Cells.Item(113, 8).Value = 4703  # H113: 3467.3 -> 4703
$ws.Cells.Item(113, 9).Value = 4111  # I113: 2939.2856 -> 4111
$ws.Cells.Item(113, 10).Value = 4999  # J113: 4699.3335 -> 4999
$ws.Cells.Item(113, 11).Value = 4111  # K113: 2939.2856 -> 4111
$ws.Cells.Item(113, 12).Value = 4999  # L113: 4699.3335 -> 4999
$ws.Cells.Item(113, 13).Value = -1941  # M113: -769.2856000000002 -> -1941
$ws.Cells.Item(113, 14).Value = -9339  # N113: -9039.333500000001 -> -9339
# Row 132
$ws.Cells.Item(132, 8).Value = 3083.1633  # H132: 3091.3264 -> 3083.1633
$ws.Cells.Item(132, 9).Value = 2056.9697  # I132: 2040.5883 -> 2056.9697
$ws.Cells.Item(132, 10).Value = 5199.6875  # J132: 5473 -> 5199.6875
$ws.Cells.Item(132, 11).Value = 6170.909100000001  # K132: 6121.7649 -> 6170.909100000001
$ws.Cells.Item(132, 12).Value = 15599.0625  # L132: 16419 -> 15599.0625
$ws.Cells.Item(132, 13).Value = -3640.909100000001  # M132: -3591.7649 -> -3640.909100000001
$ws.Cells.Item(132, 14).Value = -20659.0625  # N132: -21479 -> -20659.0625

# --- Sheet 7: LTW ---
$ws = $wb.Worksheets.Item(7)
# Row 61
$ws.Cells.Item(61, 8).Value = 51453.285  # H61: 51599.617 -> 51453.285
$ws.Cells.Item(61, 9).Value = 61942.293  # I61: 62123.06 -> 61942.293
$ws.Cells.Item(61, 11).Value = 61942.293  # K61: 62123.06 -> 61942.293
$ws.Cells.Item(61, 13).Value = -61740.293  # M61: -61921.06 -> -61740.293
# Row 99
$ws.Cells.Item(99, 8).Value = 13749.5  # H99: 13999.5 -> 13749.5
$ws.Cells.Item(99, 9).Value = 7500  # I99: 8000 -> 7500
$ws.Cells.Item(99, 11).Value = 7500  # K99: 8000 -> 7500
$ws.Cells.Item(99, 13).Value = -4505  # M99: -5005 -> -4505
# Row 113
$ws.Cells.Item(113, 8).Value = 51453.285  # H113: 51599.617 -> 51453.285
$ws.Cells.Item(113, 9).Value = 61942.293  # I113: 62123.06 -> 61942.293
$ws.Cells.Item(113, 11).Value = 61942.293  # K113: 62123.06 -> 61942.293
$ws.Cells.Item(113, 13).Value = -59772.293  # M113: -59953.06 -> -59772.293
# Row 122
$ws.Cells.Item(122, 8).Value = 462817.53  # H122: 424999.34 -> 462817.53
$ws.Cells.Item(122, 9).Value = 776884.3  # I122: 674499.6 -> 776884.3
$ws.Cells.Item(122, 11).Value = 2330652.9  # K122: 2023498.8 -> 2330652.9
$ws.Cells.Item(122, 13).Value = -2328202.9  # M122: -2021048.8 -> -2328202.9
# Row 132
$ws.Cells.Item(132, 8).Value = 5115.7915  # H132: 5290.091 -> 5115.7915
$ws.Cells.Item(132, 9).Value = 5138.769  # I132: 5491.5454 -> 5138.769
$ws.Cells.Item(132, 11).Value = 15416.307  # K132: 16474.6362 -> 15416.307
$ws.Cells.Item(132, 13).Value = -12886.307  # M132: -13944.6362 -> -12886.307
# Row 136
$ws.Cells.Item(136, 8).Value = 5299.4165  # H136: 6438.7 -> 5299.4165
$ws.Cells.Item(136, 9).Value = 4959.3  # I136: 5985.875 -> 4959.3
$ws.Cells.Item(136, 10).Value = 7000  # J136: 8250 -> 7000
$ws.Cells.Item(136, 11).Value = 14877.9  # K136: 17957.625 -> 14877.9
$ws.Cells.Item(136, 12).Value = 21000  # L136: 24750 -> 21000
$ws.Cells.Item(136, 13).Value = -12327.9  # M136: -15407.625 -> -12327.9
$ws.Cells.Item(136, 14).Value = -26100  # N136: -29850 -> -26100

# --- Sheet 8: WVR ---
$ws = $wb.Worksheets.Item(8)
# Row 9
$ws.Cells.Item(9, 8).Value = 37166  # H9: 46247.25 -> 37166
$ws.Cells.Item(9, 9).Value = 40997.8  # I9: 46247.25 -> 40997.8
$ws.Cells.Item(9, 10).Value = 18007  # J9: 0 -> 18007
$ws.Cells.Item(9, 11).Value = 40997.8  # K9: 46247.25 -> 40997.8
$ws.Cells.Item(9, 12).Value = 18007  # L9: 0 -> 18007
$ws.Cells.Item(9, 13).Value = -40857.8  # M9: -46107.25 -> -40857.8
$ws.Cells.Item(9, 14).Value = -18287  # N9: add (new -18287)
# Row 122
$ws.Cells.Item(122, 8).Value = 2982.7368  # H122: 2456.9167 -> 2982.7368
$ws.Cells.Item(122, 9).Value = 3387.4  # I122: 2411.2 -> 3387.4
$ws.Cells.Item(122, 11).Value = 10162.2  # K122: 7233.599999999999 -> 10162.2
$ws.Cells.Item(122, 13).Value = -7712.200000000001  # M122: -4783.599999999999 -> -7712.200000000001
# Row 126
$ws.Cells.Item(126, 8).Value = 14491.1  # H126: 15656.333 -> 14491.1
$ws.Cells.Item(126, 9).Value = 16488.875  # I126: 18272.428 -> 16488.875
$ws.Cells.Item(126, 11).Value = 49466.625  # K126: 54817.284 -> 49466.625
$ws.Cells.Item(126, 13).Value = -46996.625  # M126: -52347.284 -> -46996.625
# Row 132
$ws.Cells.Item(132, 8).Value = 3182.675  # H132: 3125.0244 -> 3182.675
$ws.Cells.Item(132, 9).Value = 2991.8215  # I132: 2923 -> 2991.8215
$ws.Cells.Item(132, 10).Value = 3628  # J132: 3613.25 -> 3628
$ws.Cells.Item(132, 11).Value = 8975.4645  # K132: 8769 -> 8975.4645
$ws.Cells.Item(132, 12).Value = 10884  # L132: 10839.75 -> 10884
$ws.Cells.Item(132, 13).Value = -6445.4645  # M132: -6239 -> -6445.4645
$ws.Cells.Item(132, 14).Value = -15944  # N132: -15899.75 -> -15944

